$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1985.8
$ws.Range("K5").Value = 1370.5
$ws.Range("K6").Value = 2311.9

$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)
